# "1st changes of mifos to finflux"
# The "Repayment schedule" sheet gets a new (blank) column inserted right
# before the existing "Late" column (column N), pushing "Late" -> O and
# "Outstanding" -> P -> Q. The newly selected/active sheet becomes
# "Repayment schedule" (it was "Edit Repayment Schedule" before).

$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

# Insert a brand-new blank column at N, shifting "Late" (N) -> O and
# "Outstanding" (P) -> Q.
$wsSchedule.Columns("N").Insert() | Out-Null

# The workbook now opens on the "Repayment schedule" tab, with G20 selected,
# instead of "Edit Repayment Schedule".
$wsSchedule.Activate() | Out-Null
$wsSchedule.Range("G20").Select() | Out-Null
